$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G so the old "d=7"/"d=10" columns
# (G:H) shift right to (H:I), making room for the new "d=6" column.
$ws.Range("G:G").Insert()

# New header for the inserted column (copy formatting from the
# neighbouring "d=5" header cell).
$ws.Range("G1").Value = "d=6"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for the inserted column (rows 2-6).
$ws.Range("G2").Value = 97.97140634294182
$ws.Range("G3").Value = 98.15734319214975
$ws.Range("G4").Value = 98.06677039628302
$ws.Range("G5").Value = 98.00915906376689
$ws.Range("G6").Value = 98.02920870325303
